# Apply the "Deploying to gh-pages" content refresh:
#   - bump Version / Date in the Metadata sheet
#   - fill in the Publisher value
#   - replace the duplicated "Contact" rows with a single Jurisdiction row
#   - sync the root Extension row's Short/Definition text on the Elements sheet

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value, previously blank
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 used to be "Contact" / "No display for ContactDetail"; retarget it
# to the Jurisdiction property and drop the duplicate row that followed it.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")

# The root Extension row's Short/Definition columns (K/L) get synced with
# the StructureDefinition's actual Title/Description instead of the
# generic placeholder text.
$elements.Range("K2").Value = "Legally Documented Sex"
$elements.Range("L2").Value = "Sex as defined or amended on the patient's legal documentation, for example a birth certificate or driver's license"
